$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly record was inserted right after the existing row 840, pushing
# the former rows 841-886 down to 842-887 (and growing the used range from
# A1:R886 to A1:R887).
$ws.Rows.Item(841).Insert()

# Populate the newly inserted row 841 with the new record's data.
$ws.Cells.Item(841, 1).Value = 9
$ws.Cells.Item(841, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(841, 3).Value = "Metropolitana"
$ws.Cells.Item(841, 4).Value = 45267
$ws.Cells.Item(841, 5).Value = 13
$ws.Cells.Item(841, 6).Value = 100112031
$ws.Cells.Item(841, 7).Value = "Poroto verde"
$ws.Cells.Item(841, 8).Value = "Magnum"
$ws.Cells.Item(841, 9).Value = "Primera"
$ws.Cells.Item(841, 10).Value = 70
$ws.Cells.Item(841, 11).Value = 45000
$ws.Cells.Item(841, 12).Value = 46000
$ws.Cells.Item(841, 13).Value = 45500
$ws.Cells.Item(841, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(841, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(841, 16).Value = 1820
$ws.Cells.Item(841, 17).Value = 25
$ws.Cells.Item(841, 18).Value = "Hortaliza"
